$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1065648418536256
$ws.Range("D2").Value = 0.08088279698951339
$ws.Range("E2").Value = 0.1222829297080565
$ws.Range("F2").Value = 1.567795569274352
$ws.Range("G2").Value = 0.00249042828720749
$ws.Range("I2").Value = 1.112500021953245
$ws.Range("J2").Value = 0.2044546736230046
$ws.Range("L2").Value = 0.2329032696459166
$ws.Range("M2").Value = 11.1519172752308
$ws.Range("O2").Value = 4.356235532505224

$ws.Range("C3").Value = 0.1092142457331988
$ws.Range("D3").Value = 0.08172115098177812
$ws.Range("E3").Value = 0.1196594475243344
$ws.Range("F3").Value = 1.614236888944383
$ws.Range("G3").Value = 0.002495694760960494
$ws.Range("I3").Value = 1.145267340993655
$ws.Range("J3").Value = 0.1956809204689307
$ws.Range("L3").Value = 0.2213196325890152
$ws.Range("M3").Value = 9.83780210534934
$ws.Range("O3").Value = 4.470633553575965

$ws.Range("C4").Value = 0.1109403052833144
$ws.Range("D4").Value = 0.08228556949231347
$ws.Range("E4").Value = 0.1181627814050685
$ws.Range("F4").Value = 1.644714586181884
$ws.Range("G4").Value = 0.002499069716987789
$ws.Range("I4").Value = 1.166798339124703
$ws.Range("J4").Value = 0.1904931070311875
$ws.Range("L4").Value = 0.2143896021389367
$ws.Range("M4").Value = 9.027308744481047
$ws.Range("O4").Value = 4.546372335344444

$ws.Range("C5").Value = 0.1116684650570363
$ws.Range("D5").Value = 0.08252801653854647
$ws.Range("E5").Value = 0.117581124659992
$ws.Range("F5").Value = 1.657622655888979
$ws.Range("G5").Value = 0.002500480745365305
$ws.Range("I5").Value = 1.17592348697999
$ws.Range("J5").Value = 0.1884280302143111
$ws.Range("L5").Value = 0.2116104444556299
$ws.Range("M5").Value = 8.696095277139705
$ws.Range("O5").Value = 4.578603640126815

$ws.Range("C6").Value = 0.1117908662072153
$ws.Range("D6").Value = 0.08256902483990558
$ws.Range("E6").Value = 0.1174862341672167
$ws.Range("F6").Value = 1.659795359179128
$ws.Range("G6").Value = 0.002500717207035267
$ws.Range("I6").Value = 1.177459802190388
$ws.Range("J6").Value = 0.1880880539385146
$ws.Range("L6").Value = 0.211151651506583
$ws.Range("M6").Value = 8.641040799143184
$ws.Range("O6").Value = 4.584037747539611

$ws.Range("C7").Value = 0.1109500254361748
$ws.Range("D7").Value = 0.08228878890032831
$ws.Range("E7").Value = 0.1181548231939047
$ws.Range("F7").Value = 1.644886699397063
$ws.Range("G7").Value = 0.002499088601894707
$ws.Range("I7").Value = 1.166919987633253
$ws.Range("J7").Value = 0.1904650597102773
$ws.Range("L7").Value = 0.214351940932886
$ws.Range("M7").Value = 9.022845676672375
$ws.Range("O7").Value = 4.546801502714871

$ws.Range("C8").Value = 0.1074576237473366
$ws.Range("D8").Value = 0.08116152401294841
$ws.Range("E8").Value = 0.1213543710674472
$ws.Range("F8").Value = 1.583397999196521
$ws.Range("G8").Value = 0.002492214940769621
$ws.Range("I8").Value = 1.123502840451874
$ws.Range("J8").Value = 0.2013874188810121
$ws.Range("L8").Value = 0.2288707676887327
$ws.Range("M8").Value = 10.69954818896065
$ws.Range("O8").Value = 4.394529619728047

$ws.Range("C9").Value = 0.1014060297364097
$ws.Range("D9").Value = 0.07934733701458896
$ws.Range("E9").Value = 0.1285572857188697
$ws.Range("F9").Value = 1.478628955418898
$ws.Range("G9").Value = 0.002479848938127008
$ws.Range("I9").Value = 1.049737654992072
$ws.Range("J9").Value = 0.2244429513380197
$ws.Range("L9").Value = 0.2588386245368355
$ws.Range("M9").Value = 13.96009794403943
$ws.Range("O9").Value = 4.140236624260439

$ws.Range("C10").Value = 0.09745838865647016
$ws.Range("D10").Value = 0.07825956124137434
$ws.Range("E10").Value = 0.1344499326507602
$ws.Range("F10").Value = 1.411620463770113
$ws.Range("G10").Value = 0.002471430499831283
$ws.Range("I10").Value = 1.002716676906161
$ws.Range("J10").Value = 0.2424638378358992
$ws.Range("L10").Value = 0.2818450863568955
$ws.Range("M10").Value = 16.34117801188609
$ws.Range("O10").Value = 3.981368686833065

$ws.Range("C11").Value = 0.0957735925526535
$ws.Range("D11").Value = 0.0778187654510063
$ws.Range("E11").Value = 0.1372693292923302
$ws.Range("F11").Value = 1.383371948583253
$ws.Range("G11").Value = 0.002467742891406368
$ws.Range("I11").Value = 0.9829357951622839
$ws.Range("J11").Value = 0.2509170144259372
$ws.Range("L11").Value = 0.2925446841184112
$ws.Range("M11").Value = 17.42185699781629
$ws.Range("O11").Value = 3.915374336634159

$ws.Range("C12").Value = 0.09515181833087638
$ws.Range("D12").Value = 0.07765969440834652
$ws.Range("E12").Value = 0.138357638597796
$ws.Range("F12").Value = 1.373002309879944
$ws.Range("G12").Value = 0.002466366694807382
$ws.Range("I12").Value = 0.975681148860879
$ws.Range("J12").Value = 0.2541564748574245
$ws.Range("L12").Value = 0.2966315802286488
$ws.Range("M12").Value = 17.83077520561966
$ws.Range("O12").Value = 3.891303988606722

$ws.Range("C13").Value = 0.09528500285256314
$ws.Range("D13").Value = 0.07769360283009519
$ws.Range("E13").Value = 0.138122320941541
$ws.Range("F13").Value = 1.375220933038584
$ws.Range("G13").Value = 0.002466662187115253
$ws.Range("I13").Value = 0.9772329998493632
$ws.Range("J13").Value = 0.2534570618847312
$ws.Range("L13").Value = 0.2957498028667374
$ws.Range("M13").Value = 17.74272041954873
$ws.Range("O13").Value = 3.896446751050405

$ws.Range("C14").Value = 0.09572211219398241
$ws.Range("D14").Value = 0.07780552081581504
$ws.Range("E14").Value = 0.1373584468039084
$ws.Range("F14").Value = 1.382512233375827
$ws.Range("G14").Value = 0.002467629266979211
$ws.Range("I14").Value = 0.9823341960668088
$ws.Range("J14").Value = 0.2511827465582996
$ws.Range("L14").Value = 0.2928802009102753
$ws.Range("M14").Value = 17.45550490565125
$ws.Range("O14").Value = 3.913375507363583

$ws.Range("C15").Value = 0.09599197465176701
$ws.Range("D15").Value = 0.0778750983558254
$ws.Range("E15").Value = 0.136893264700376
$ws.Range("F15").Value = 1.387021196793036
$ws.Range("G15").Value = 0.002468224258456993
$ws.Range("I15").Value = 0.9854896891631739
$ws.Range("J15").Value = 0.2497947192111667
$ws.Range("L15").Value = 0.2911271179572168
$ws.Range("M15").Value = 17.27953799913621
$ws.Range("O15").Value = 3.923865242838247

$ws.Range("C16").Value = 0.09757074968912605
$ws.Range("D16").Value = 0.07828946223622069
$ws.Range("E16").Value = 0.134268530823789
$ws.Range("F16").Value = 1.413512053116015
$ws.Range("G16").Value = 0.002471674334171411
$ws.Range("I16").Value = 1.00404217182323
$ws.Range("J16").Value = 0.2419166849992536
$ws.Range("L16").Value = 0.2811506793044032
$ws.Range("M16").Value = 16.27050689252047
$ws.Range("O16").Value = 3.985809372349621

$ws.Range("C17").Value = 0.0985678895932196
$ws.Range("D17").Value = 0.07855755560130717
$ws.Range("E17").Value = 0.1326943805195668
$ws.Range("F17").Value = 1.430339902002281
$ws.Range("G17").Value = 0.002473827070450848
$ws.Range("I17").Value = 1.01583888128971
$ws.Range("J17").Value = 0.2371503548845055
$ws.Range("L17").Value = 0.2750914297482581
$ws.Range("M17").Value = 15.65089340420388
$ws.Range("O17").Value = 4.025430010393677

$ws.Range("C18").Value = 0.09915185042208208
$ws.Range("D18").Value = 0.07871683723231371
$ws.Range("E18").Value = 0.1318019907956227
$ws.Range("F18").Value = 1.440228730857555
$ws.Range("G18").Value = 0.002475078643562933
$ws.Range("I18").Value = 1.022775225607703
$ws.Range("J18").Value = 0.2344328178469368
$ws.Range("L18").Value = 0.2716282200259741
$ws.Range("M18").Value = 15.29427103283427
$ws.Range("O18").Value = 4.0488083441181

$ws.Range("C19").Value = 0.09935135354265867
$ws.Range("D19").Value = 0.07877163769551743
$ws.Range("E19").Value = 0.1315020583451627
$ws.Range("F19").Value = 1.443612780344992
$ws.Range("G19").Value = 0.002475504708101421
$ws.Range("I19").Value = 1.025149582248659
$ws.Range("J19").Value = 0.2335167663045183
$ws.Range("L19").Value = 0.2704593538953475
$ws.Range("M19").Value = 15.17348278445161
$ws.Range("O19").Value = 4.056824569110233

$ws.Range("C20").Value = 0.09846066079640536
$ws.Range("D20").Value = 0.07852849014928154
$ws.Range("E20").Value = 0.1328605990243901
$ws.Range("F20").Value = 1.428526780719515
$ws.Range("G20").Value = 0.002473596525136091
$ws.Range("I20").Value = 1.014567421281576
$ws.Range("J20").Value = 0.2376552502004614
$ws.Range("L20").Value = 0.2757341685747434
$ws.Range("M20").Value = 15.71687653681835
$ws.Range("O20").Value = 4.021151176970335

$ws.Range("C21").Value = 0.0955932802250814
$ws.Range("D21").Value = 0.07777243412872537
$ws.Range("E21").Value = 0.1375822484751197
$ws.Range("F21").Value = 1.380361664429259
$ws.Range("G21").Value = 0.002467344664834697
$ws.Range("I21").Value = 0.9808294110050397
$ws.Range("J21").Value = 0.2518497111528575
$ws.Range("L21").Value = 0.2937221045032175
$ws.Range("M21").Value = 17.53987512595756
$ws.Range("O21").Value = 3.908378001100942

$ws.Range("C22").Value = 0.09381395963564643
$ws.Range("D22").Value = 0.0773240963387849
$ws.Range("E22").Value = 0.1407889291071598
$ws.Range("F22").Value = 1.350794634706808
$ws.Range("G22").Value = 0.002463376476993439
$ws.Range("I22").Value = 0.9601571023911148
$ws.Range("J22").Value = 0.2613514516842912
$ws.Range("L22").Value = 0.3056841809273294
$ws.Range("M22").Value = 18.72951951418838
$ws.Range("O22").Value = 3.840046832958535

$ws.Range("C23").Value = 0.09475486536881306
$ws.Range("D23").Value = 0.07755916585448119
$ws.Range("E23").Value = 0.139066164826545
$ws.Range("F23").Value = 1.366398115004849
$ws.Range("G23").Value = 0.002465483664677844
$ws.Range("I23").Value = 0.9710627449396938
$ws.Range("J23").Value = 0.2562590317291438
$ws.Range("L23").Value = 0.2992804091813355
$ws.Range("M23").Value = 18.09473124515767
$ws.Range("O23").Value = 3.876018879193538

$ws.Range("C24").Value = 0.09850910565269544
$ws.Range("D24").Value = 0.07854161460127074
$ws.Range("E24").Value = 0.1327854124130923
$ws.Range("F24").Value = 1.429345825904186
$ws.Range("G24").Value = 0.002473700711176031
$ws.Range("I24").Value = 1.015141768212388
$ws.Range("J24").Value = 0.2374269164824909
$ws.Range("L24").Value = 0.2754435231337737
$ws.Range("M24").Value = 15.68704680666951
$ws.Range("O24").Value = 4.023083770352628

$ws.Range("C25").Value = 0.1029565490084252
$ws.Range("D25").Value = 0.0797953461983667
$ws.Range("E25").Value = 0.1265058624042865
$ws.Range("F25").Value = 1.505245650734317
$ws.Range("G25").Value = 0.002483076290705644
$ws.Range("I25").Value = 1.068450708353964
$ws.Range("J25").Value = 0.2180222315185887
$ws.Range("L25").Value = 0.2505637354316974
$ws.Range("M25").Value = 13.08075262910694
$ws.Range("O25").Value = 4.204188609592194
